# Atualização automática da planilha
# Target sheet: "Budget" (tabColor FFDC2626), holding the budget line items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budget")
$ws.Activate()

# Delete row 5 (the "Consultoria & Implantação / SoftExpert / P03 / 80000" row)
# — rows below shift up one position.
$ws.Rows("5:5").Delete()

# Zero out the "Orçado (R$)" column (D) for the remaining budget line items.
$ws.Range("D3:D9").Value = 0

# Restore the reported selection from the edit.
$ws.Range("C12:D14").Select()
$ws.Application.ActiveWindow.RangeSelection.Item(1).Activate()
